$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in task name: "Metode- og teknolgivalg" -> "Metode- og teknologivalg"
$ws.Range("C10").Value = "Metode- og teknologivalg"

# The Gantt bar for "Metode- og teknologivalg" (row 10) shifts one day later:
# drop the colored marker on N10 and add it on R10 (copy O10's fill/format).
$ws.Range("O10").Copy()
$ws.Range("R10").PasteSpecial(-4122)
$ws.Range("N10").Clear()

# The Gantt bar for "Teknologiafsnit" (row 17) shifts one day earlier:
# add the colored marker on N17 (copy O17's fill/format) and clear it from R17.
$ws.Range("O17").Copy()
$ws.Range("N17").PasteSpecial(-4122)
$ws.Range("S17").Copy()
$ws.Range("R17").PasteSpecial(-4122)

# Column C needs to widen slightly to fit the corrected (longer) text.
$ws.Columns("C:C").ColumnWidth = 23.33

# Update the active selection left at the end of the editing session.
[void]$ws.Range("U15").Select()
